$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.608.61"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "2.094.90"
$ws.Range("E3").Value = "  +3.02%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "2.409.43"
$ws.Range("E12").Value = "  +3.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.783"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.04%  "

$ws.Range("D17").Value = "2.104.85"
$ws.Range("E17").Value = "  +3.67%  "

$ws.Range("D18").Value = "38.555.24"
$ws.Range("E18").Value = "  +2.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "

$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("E28").Value = "  +2.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.80%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("E32").Value = "  +3.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.15%  "

$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.80%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.82%  "

$ws.Range("D41").Value = "1.543.33"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.15%  "

$ws.Range("E47").Value = "  +8.40%  "

$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").Value = "2.293.93"
$ws.Range("E51").Value = "  +3.18%  "
